$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6,1).Value = 3
$ws.Cells.Item(6,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(6,3).Value = 'Coquimbo'
$ws.Cells.Item(6,4).Value = 45050
$ws.Cells.Item(6,5).Value = 5
$ws.Cells.Item(6,6).Value = 'Fruta'
$ws.Cells.Item(6,7).Value = 100104
$ws.Cells.Item(6,8).Value = 'Frutos de pepita'
$ws.Cells.Item(6,9).Value = 100104001
$ws.Cells.Item(6,10).Value = 'Granada'
$ws.Cells.Item(6,11).Value = 'Wonderfull'
$ws.Cells.Item(6,12).Value = 'Especial'
$ws.Cells.Item(6,13).Value = 56
$ws.Cells.Item(6,14).Value = 14000
$ws.Cells.Item(6,15).Value = 14000
$ws.Cells.Item(6,16).Value = 14000
$ws.Cells.Item(6,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(6,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(6,19).Value = 1000
$ws.Cells.Item(6,20).Value = 14
$ws.Cells.Item(7,1).Value = 3
$ws.Cells.Item(7,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(7,3).Value = 'Coquimbo'
$ws.Cells.Item(7,4).Value = 45050
$ws.Cells.Item(7,5).Value = 5
$ws.Cells.Item(7,6).Value = 'Fruta'
$ws.Cells.Item(7,7).Value = 100104
$ws.Cells.Item(7,8).Value = 'Frutos de pepita'
$ws.Cells.Item(7,9).Value = 100104001
$ws.Cells.Item(7,10).Value = 'Granada'
$ws.Cells.Item(7,11).Value = 'Wonderfull'
$ws.Cells.Item(7,12).Value = 'Primera'
$ws.Cells.Item(7,13).Value = 50
$ws.Cells.Item(7,14).Value = 12000
$ws.Cells.Item(7,15).Value = 12000
$ws.Cells.Item(7,16).Value = 12000
$ws.Cells.Item(7,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(7,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(7,19).Value = 857
$ws.Cells.Item(7,20).Value = 14
$ws.Cells.Item(8,1).Value = 3
$ws.Cells.Item(8,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(8,3).Value = 'Coquimbo'
$ws.Cells.Item(8,4).Value = 44252
$ws.Cells.Item(8,5).Value = 5
$ws.Cells.Item(8,6).Value = 'Fruta'
$ws.Cells.Item(8,7).Value = 100104
$ws.Cells.Item(8,8).Value = 'Frutos de pepita'
$ws.Cells.Item(8,9).Value = 100104001
$ws.Cells.Item(8,10).Value = 'Granada'
$ws.Cells.Item(8,11).Value = 'Wonderfull'
$ws.Cells.Item(8,12).Value = 'Primera'
$ws.Cells.Item(8,13).Value = 60
$ws.Cells.Item(8,14).Value = 14000
$ws.Cells.Item(8,15).Value = 14000
$ws.Cells.Item(8,16).Value = 14000
$ws.Cells.Item(8,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(8,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(8,19).Value = 1000
$ws.Cells.Item(8,20).Value = 14
$ws.Cells.Item(9,1).Value = 3
$ws.Cells.Item(9,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(9,3).Value = 'Coquimbo'
$ws.Cells.Item(9,4).Value = 44614
$ws.Cells.Item(9,5).Value = 5
$ws.Cells.Item(9,6).Value = 'Fruta'
$ws.Cells.Item(9,7).Value = 100104
$ws.Cells.Item(9,8).Value = 'Frutos de pepita'
$ws.Cells.Item(9,9).Value = 100104001
$ws.Cells.Item(9,10).Value = 'Granada'
$ws.Cells.Item(9,11).Value = 'Wonderfull'
$ws.Cells.Item(9,12).Value = 'Primera'
$ws.Cells.Item(9,13).Value = 54
$ws.Cells.Item(9,14).Value = 14000
$ws.Cells.Item(9,15).Value = 14000
$ws.Cells.Item(9,16).Value = 14000
$ws.Cells.Item(9,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(9,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(9,19).Value = 1000
$ws.Cells.Item(9,20).Value = 14
$ws.Cells.Item(10,1).Value = 3
$ws.Cells.Item(10,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(10,3).Value = 'Coquimbo'
$ws.Cells.Item(10,4).Value = 44245
$ws.Cells.Item(10,5).Value = 5
$ws.Cells.Item(10,6).Value = 'Fruta'
$ws.Cells.Item(10,7).Value = 100104
$ws.Cells.Item(10,8).Value = 'Frutos de pepita'
$ws.Cells.Item(10,9).Value = 100104001
$ws.Cells.Item(10,10).Value = 'Granada'
$ws.Cells.Item(10,11).Value = 'Wonderfull'
$ws.Cells.Item(10,12).Value = 'Primera'
$ws.Cells.Item(10,13).Value = 50
$ws.Cells.Item(10,14).Value = 15000
$ws.Cells.Item(10,15).Value = 15000
$ws.Cells.Item(10,16).Value = 15000
$ws.Cells.Item(10,17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(10,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(10,19).Value = 1000
$ws.Cells.Item(10,20).Value = 15
$ws.Cells.Item(11,1).Value = 3
$ws.Cells.Item(11,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(11,3).Value = 'Coquimbo'
$ws.Cells.Item(11,4).Value = 44627
$ws.Cells.Item(11,5).Value = 5
$ws.Cells.Item(11,6).Value = 'Fruta'
$ws.Cells.Item(11,7).Value = 100104
$ws.Cells.Item(11,8).Value = 'Frutos de pepita'
$ws.Cells.Item(11,9).Value = 100104001
$ws.Cells.Item(11,10).Value = 'Granada'
$ws.Cells.Item(11,11).Value = 'Wonderfull'
$ws.Cells.Item(11,12).Value = 'Primera'
$ws.Cells.Item(11,13).Value = 56
$ws.Cells.Item(11,14).Value = 17000
$ws.Cells.Item(11,15).Value = 17000
$ws.Cells.Item(11,16).Value = 17000
$ws.Cells.Item(11,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(11,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(11,19).Value = 1214
$ws.Cells.Item(11,20).Value = 14
$ws.Cells.Item(12,1).Value = 3
$ws.Cells.Item(12,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(12,3).Value = 'Coquimbo'
$ws.Cells.Item(12,4).Value = 45001
$ws.Cells.Item(12,5).Value = 5
$ws.Cells.Item(12,6).Value = 'Fruta'
$ws.Cells.Item(12,7).Value = 100104
$ws.Cells.Item(12,8).Value = 'Frutos de pepita'
$ws.Cells.Item(12,9).Value = 100104001
$ws.Cells.Item(12,10).Value = 'Granada'
$ws.Cells.Item(12,11).Value = 'Wonderfull'
$ws.Cells.Item(12,12).Value = 'Primera'
$ws.Cells.Item(12,13).Value = 50
$ws.Cells.Item(12,14).Value = 16000
$ws.Cells.Item(12,15).Value = 16000
$ws.Cells.Item(12,16).Value = 16000
$ws.Cells.Item(12,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(12,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(12,19).Value = 1143
$ws.Cells.Item(12,20).Value = 14
$ws.Cells.Item(13,1).Value = 3
$ws.Cells.Item(13,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(13,3).Value = 'Coquimbo'
$ws.Cells.Item(13,4).Value = 44320
$ws.Cells.Item(13,5).Value = 5
$ws.Cells.Item(13,6).Value = 'Fruta'
$ws.Cells.Item(13,7).Value = 100104
$ws.Cells.Item(13,8).Value = 'Frutos de pepita'
$ws.Cells.Item(13,9).Value = 100104001
$ws.Cells.Item(13,10).Value = 'Granada'
$ws.Cells.Item(13,11).Value = 'Wonderfull'
$ws.Cells.Item(13,12).Value = 'Primera'
$ws.Cells.Item(13,13).Value = 45
$ws.Cells.Item(13,14).Value = 14000
$ws.Cells.Item(13,15).Value = 14000
$ws.Cells.Item(13,16).Value = 14000
$ws.Cells.Item(13,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(13,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(13,19).Value = 1000
$ws.Cells.Item(13,20).Value = 14
$ws.Cells.Item(14,1).Value = 3
$ws.Cells.Item(14,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(14,3).Value = 'Coquimbo'
$ws.Cells.Item(14,4).Value = 45014
$ws.Cells.Item(14,5).Value = 5
$ws.Cells.Item(14,6).Value = 'Fruta'
$ws.Cells.Item(14,7).Value = 100104
$ws.Cells.Item(14,8).Value = 'Frutos de pepita'
$ws.Cells.Item(14,9).Value = 100104001
$ws.Cells.Item(14,10).Value = 'Granada'
$ws.Cells.Item(14,11).Value = 'Wonderfull'
$ws.Cells.Item(14,12).Value = 'Primera'
$ws.Cells.Item(14,13).Value = 60
$ws.Cells.Item(14,14).Value = 15000
$ws.Cells.Item(14,15).Value = 15000
$ws.Cells.Item(14,16).Value = 15000
$ws.Cells.Item(14,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(14,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(14,19).Value = 1071
$ws.Cells.Item(14,20).Value = 14
$ws.Cells.Item(15,1).Value = 3
$ws.Cells.Item(15,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(15,3).Value = 'Coquimbo'
$ws.Cells.Item(15,4).Value = 44260
$ws.Cells.Item(15,5).Value = 5
$ws.Cells.Item(15,6).Value = 'Fruta'
$ws.Cells.Item(15,7).Value = 100104
$ws.Cells.Item(15,8).Value = 'Frutos de pepita'
$ws.Cells.Item(15,9).Value = 100104001
$ws.Cells.Item(15,10).Value = 'Granada'
$ws.Cells.Item(15,11).Value = 'Wonderfull'
$ws.Cells.Item(15,12).Value = 'Primera'
$ws.Cells.Item(15,13).Value = 56
$ws.Cells.Item(15,14).Value = 13000
$ws.Cells.Item(15,15).Value = 13000
$ws.Cells.Item(15,16).Value = 13000
$ws.Cells.Item(15,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(15,18).Value = 'Provincia del Elquí'
$ws.Cells.Item(15,19).Value = 929
$ws.Cells.Item(15,20).Value = 14
$ws.Cells.Item(16,1).Value = 3
$ws.Cells.Item(16,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(16,3).Value = 'Coquimbo'
$ws.Cells.Item(16,4).Value = 44588
$ws.Cells.Item(16,5).Value = 5
$ws.Cells.Item(16,6).Value = 'Fruta'
$ws.Cells.Item(16,7).Value = 100104
$ws.Cells.Item(16,8).Value = 'Frutos de pepita'
$ws.Cells.Item(16,9).Value = 100104001
$ws.Cells.Item(16,10).Value = 'Granada'
$ws.Cells.Item(16,11).Value = 'Wonderfull'
$ws.Cells.Item(16,12).Value = 'Primera'
$ws.Cells.Item(16,13).Value = 85
$ws.Cells.Item(16,14).Value = 19000
$ws.Cells.Item(16,15).Value = 20000
$ws.Cells.Item(16,16).Value = 19529
$ws.Cells.Item(16,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(16,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(16,19).Value = 1395
$ws.Cells.Item(16,20).Value = 14
$ws.Cells.Item(17,1).Value = 3
$ws.Cells.Item(17,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(17,3).Value = 'Coquimbo'
$ws.Cells.Item(17,4).Value = 45044
$ws.Cells.Item(17,5).Value = 5
$ws.Cells.Item(17,6).Value = 'Fruta'
$ws.Cells.Item(17,7).Value = 100104
$ws.Cells.Item(17,8).Value = 'Frutos de pepita'
$ws.Cells.Item(17,9).Value = 100104001
$ws.Cells.Item(17,10).Value = 'Granada'
$ws.Cells.Item(17,11).Value = 'Wonderfull'
$ws.Cells.Item(17,12).Value = 'Especial'
$ws.Cells.Item(17,13).Value = 30
$ws.Cells.Item(17,14).Value = 16000
$ws.Cells.Item(17,15).Value = 16000
$ws.Cells.Item(17,16).Value = 16000
$ws.Cells.Item(17,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(17,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(17,19).Value = 1143
$ws.Cells.Item(17,20).Value = 14
$ws.Cells.Item(18,1).Value = 3
$ws.Cells.Item(18,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(18,3).Value = 'Coquimbo'
$ws.Cells.Item(18,4).Value = 45044
$ws.Cells.Item(18,5).Value = 5
$ws.Cells.Item(18,6).Value = 'Fruta'
$ws.Cells.Item(18,7).Value = 100104
$ws.Cells.Item(18,8).Value = 'Frutos de pepita'
$ws.Cells.Item(18,9).Value = 100104001
$ws.Cells.Item(18,10).Value = 'Granada'
$ws.Cells.Item(18,11).Value = 'Wonderfull'
$ws.Cells.Item(18,12).Value = 'Primera'
$ws.Cells.Item(18,13).Value = 30
$ws.Cells.Item(18,14).Value = 14000
$ws.Cells.Item(18,15).Value = 14000
$ws.Cells.Item(18,16).Value = 14000
$ws.Cells.Item(18,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(18,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(18,19).Value = 1000
$ws.Cells.Item(18,20).Value = 14
$ws.Cells.Item(19,1).Value = 3
$ws.Cells.Item(19,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(19,3).Value = 'Coquimbo'
$ws.Cells.Item(19,4).Value = 45015
$ws.Cells.Item(19,5).Value = 5
$ws.Cells.Item(19,6).Value = 'Fruta'
$ws.Cells.Item(19,7).Value = 100104
$ws.Cells.Item(19,8).Value = 'Frutos de pepita'
$ws.Cells.Item(19,9).Value = 100104001
$ws.Cells.Item(19,10).Value = 'Granada'
$ws.Cells.Item(19,11).Value = 'Wonderfull'
$ws.Cells.Item(19,12).Value = 'Primera'
$ws.Cells.Item(19,13).Value = 56
$ws.Cells.Item(19,14).Value = 15000
$ws.Cells.Item(19,15).Value = 15000
$ws.Cells.Item(19,16).Value = 15000
$ws.Cells.Item(19,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(19,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(19,19).Value = 1071
$ws.Cells.Item(19,20).Value = 14
$ws.Cells.Item(20,1).Value = 3
$ws.Cells.Item(20,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(20,3).Value = 'Coquimbo'
$ws.Cells.Item(20,4).Value = 44312
$ws.Cells.Item(20,5).Value = 5
$ws.Cells.Item(20,6).Value = 'Fruta'
$ws.Cells.Item(20,7).Value = 100104
$ws.Cells.Item(20,8).Value = 'Frutos de pepita'
$ws.Cells.Item(20,9).Value = 100104001
$ws.Cells.Item(20,10).Value = 'Granada'
$ws.Cells.Item(20,11).Value = 'Wonderfull'
$ws.Cells.Item(20,12).Value = 'Primera'
$ws.Cells.Item(20,13).Value = 68
$ws.Cells.Item(20,14).Value = 14000
$ws.Cells.Item(20,15).Value = 14000
$ws.Cells.Item(20,16).Value = 14000
$ws.Cells.Item(20,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(20,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(20,19).Value = 1000
$ws.Cells.Item(20,20).Value = 14
$ws.Cells.Item(21,1).Value = 3
$ws.Cells.Item(21,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(21,3).Value = 'Coquimbo'
$ws.Cells.Item(21,4).Value = 44313
$ws.Cells.Item(21,5).Value = 5
$ws.Cells.Item(21,6).Value = 'Fruta'
$ws.Cells.Item(21,7).Value = 100104
$ws.Cells.Item(21,8).Value = 'Frutos de pepita'
$ws.Cells.Item(21,9).Value = 100104001
$ws.Cells.Item(21,10).Value = 'Granada'
$ws.Cells.Item(21,11).Value = 'Wonderfull'
$ws.Cells.Item(21,12).Value = 'Primera'
$ws.Cells.Item(21,13).Value = 36
$ws.Cells.Item(21,14).Value = 14000
$ws.Cells.Item(21,15).Value = 14000
$ws.Cells.Item(21,16).Value = 14000
$ws.Cells.Item(21,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(21,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(21,19).Value = 1000
$ws.Cells.Item(21,20).Value = 14
$ws.Cells.Item(22,1).Value = 3
$ws.Cells.Item(22,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(22,3).Value = 'Coquimbo'
$ws.Cells.Item(22,4).Value = 44270
$ws.Cells.Item(22,5).Value = 5
$ws.Cells.Item(22,6).Value = 'Fruta'
$ws.Cells.Item(22,7).Value = 100104
$ws.Cells.Item(22,8).Value = 'Frutos de pepita'
$ws.Cells.Item(22,9).Value = 100104001
$ws.Cells.Item(22,10).Value = 'Granada'
$ws.Cells.Item(22,11).Value = 'Wonderfull'
$ws.Cells.Item(22,12).Value = 'Primera'
$ws.Cells.Item(22,13).Value = 85
$ws.Cells.Item(22,14).Value = 12000
$ws.Cells.Item(22,15).Value = 12000
$ws.Cells.Item(22,16).Value = 12000
$ws.Cells.Item(22,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(22,18).Value = 'Provincia del Elquí'
$ws.Cells.Item(22,19).Value = 857
$ws.Cells.Item(22,20).Value = 14
$ws.Cells.Item(23,1).Value = 3
$ws.Cells.Item(23,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(23,3).Value = 'Coquimbo'
$ws.Cells.Item(23,4).Value = 44239
$ws.Cells.Item(23,5).Value = 5
$ws.Cells.Item(23,6).Value = 'Fruta'
$ws.Cells.Item(23,7).Value = 100104
$ws.Cells.Item(23,8).Value = 'Frutos de pepita'
$ws.Cells.Item(23,9).Value = 100104001
$ws.Cells.Item(23,10).Value = 'Granada'
$ws.Cells.Item(23,11).Value = 'Wonderfull'
$ws.Cells.Item(23,12).Value = 'Primera'
$ws.Cells.Item(23,13).Value = 70
$ws.Cells.Item(23,14).Value = 15000
$ws.Cells.Item(23,15).Value = 15000
$ws.Cells.Item(23,16).Value = 15000
$ws.Cells.Item(23,17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(23,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(23,19).Value = 1000
$ws.Cells.Item(23,20).Value = 15
$ws.Cells.Item(24,1).Value = 3
$ws.Cells.Item(24,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(24,3).Value = 'Coquimbo'
$ws.Cells.Item(24,4).Value = 44259
$ws.Cells.Item(24,5).Value = 5
$ws.Cells.Item(24,6).Value = 'Fruta'
$ws.Cells.Item(24,7).Value = 100104
$ws.Cells.Item(24,8).Value = 'Frutos de pepita'
$ws.Cells.Item(24,9).Value = 100104001
$ws.Cells.Item(24,10).Value = 'Granada'
$ws.Cells.Item(24,11).Value = 'Wonderfull'
$ws.Cells.Item(24,12).Value = 'Primera'
$ws.Cells.Item(24,13).Value = 80
$ws.Cells.Item(24,14).Value = 12000
$ws.Cells.Item(24,15).Value = 12000
$ws.Cells.Item(24,16).Value = 12000
$ws.Cells.Item(24,17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(24,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(24,19).Value = 800
$ws.Cells.Item(24,20).Value = 15
$ws.Cells.Item(25,1).Value = 3
$ws.Cells.Item(25,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(25,3).Value = 'Coquimbo'
$ws.Cells.Item(25,4).Value = 44278
$ws.Cells.Item(25,5).Value = 5
$ws.Cells.Item(25,6).Value = 'Fruta'
$ws.Cells.Item(25,7).Value = 100104
$ws.Cells.Item(25,8).Value = 'Frutos de pepita'
$ws.Cells.Item(25,9).Value = 100104001
$ws.Cells.Item(25,10).Value = 'Granada'
$ws.Cells.Item(25,11).Value = 'Wonderfull'
$ws.Cells.Item(25,12).Value = 'Primera'
$ws.Cells.Item(25,13).Value = 45
$ws.Cells.Item(25,14).Value = 13000
$ws.Cells.Item(25,15).Value = 13000
$ws.Cells.Item(25,16).Value = 13000
$ws.Cells.Item(25,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(25,18).Value = 'Provincia del Elquí'
$ws.Cells.Item(25,19).Value = 929
$ws.Cells.Item(25,20).Value = 14
$ws.Cells.Item(26,1).Value = 3
$ws.Cells.Item(26,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(26,3).Value = 'Coquimbo'
$ws.Cells.Item(26,4).Value = 45042
$ws.Cells.Item(26,5).Value = 5
$ws.Cells.Item(26,6).Value = 'Fruta'
$ws.Cells.Item(26,7).Value = 100104
$ws.Cells.Item(26,8).Value = 'Frutos de pepita'
$ws.Cells.Item(26,9).Value = 100104001
$ws.Cells.Item(26,10).Value = 'Granada'
$ws.Cells.Item(26,11).Value = 'Wonderfull'
$ws.Cells.Item(26,12).Value = 'Especial'
$ws.Cells.Item(26,13).Value = 50
$ws.Cells.Item(26,14).Value = 17000
$ws.Cells.Item(26,15).Value = 17000
$ws.Cells.Item(26,16).Value = 17000
$ws.Cells.Item(26,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(26,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(26,19).Value = 1214
$ws.Cells.Item(26,20).Value = 14
$ws.Cells.Item(27,1).Value = 3
$ws.Cells.Item(27,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(27,3).Value = 'Coquimbo'
$ws.Cells.Item(27,4).Value = 45042
$ws.Cells.Item(27,5).Value = 5
$ws.Cells.Item(27,6).Value = 'Fruta'
$ws.Cells.Item(27,7).Value = 100104
$ws.Cells.Item(27,8).Value = 'Frutos de pepita'
$ws.Cells.Item(27,9).Value = 100104001
$ws.Cells.Item(27,10).Value = 'Granada'
$ws.Cells.Item(27,11).Value = 'Wonderfull'
$ws.Cells.Item(27,12).Value = 'Primera'
$ws.Cells.Item(27,13).Value = 50
$ws.Cells.Item(27,14).Value = 14000
$ws.Cells.Item(27,15).Value = 14000
$ws.Cells.Item(27,16).Value = 14000
$ws.Cells.Item(27,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(27,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(27,19).Value = 1000
$ws.Cells.Item(27,20).Value = 14
$ws.Cells.Item(28,1).Value = 3
$ws.Cells.Item(28,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(28,3).Value = 'Coquimbo'
$ws.Cells.Item(28,4).Value = 44314
$ws.Cells.Item(28,5).Value = 5
$ws.Cells.Item(28,6).Value = 'Fruta'
$ws.Cells.Item(28,7).Value = 100104
$ws.Cells.Item(28,8).Value = 'Frutos de pepita'
$ws.Cells.Item(28,9).Value = 100104001
$ws.Cells.Item(28,10).Value = 'Granada'
$ws.Cells.Item(28,11).Value = 'Wonderfull'
$ws.Cells.Item(28,12).Value = 'Primera'
$ws.Cells.Item(28,13).Value = 56
$ws.Cells.Item(28,14).Value = 14000
$ws.Cells.Item(28,15).Value = 14000
$ws.Cells.Item(28,16).Value = 14000
$ws.Cells.Item(28,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(28,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(28,19).Value = 1000
$ws.Cells.Item(28,20).Value = 14
$ws.Cells.Item(29,1).Value = 3
$ws.Cells.Item(29,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(29,3).Value = 'Coquimbo'
$ws.Cells.Item(29,4).Value = 44316
$ws.Cells.Item(29,5).Value = 5
$ws.Cells.Item(29,6).Value = 'Fruta'
$ws.Cells.Item(29,7).Value = 100104
$ws.Cells.Item(29,8).Value = 'Frutos de pepita'
$ws.Cells.Item(29,9).Value = 100104001
$ws.Cells.Item(29,10).Value = 'Granada'
$ws.Cells.Item(29,11).Value = 'Wonderfull'
$ws.Cells.Item(29,12).Value = 'Primera'
$ws.Cells.Item(29,13).Value = 48
$ws.Cells.Item(29,14).Value = 14000
$ws.Cells.Item(29,15).Value = 14000
$ws.Cells.Item(29,16).Value = 14000
$ws.Cells.Item(29,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(29,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(29,19).Value = 1000
$ws.Cells.Item(29,20).Value = 14
$ws.Cells.Item(30,1).Value = 3
$ws.Cells.Item(30,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(30,3).Value = 'Coquimbo'
$ws.Cells.Item(30,4).Value = 45006
$ws.Cells.Item(30,5).Value = 5
$ws.Cells.Item(30,6).Value = 'Fruta'
$ws.Cells.Item(30,7).Value = 100104
$ws.Cells.Item(30,8).Value = 'Frutos de pepita'
$ws.Cells.Item(30,9).Value = 100104001
$ws.Cells.Item(30,10).Value = 'Granada'
$ws.Cells.Item(30,11).Value = 'Wonderfull'
$ws.Cells.Item(30,12).Value = 'Primera'
$ws.Cells.Item(30,13).Value = 40
$ws.Cells.Item(30,14).Value = 16000
$ws.Cells.Item(30,15).Value = 16000
$ws.Cells.Item(30,16).Value = 16000
$ws.Cells.Item(30,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(30,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(30,19).Value = 1143
$ws.Cells.Item(30,20).Value = 14
$ws.Cells.Item(31,1).Value = 3
$ws.Cells.Item(31,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(31,3).Value = 'Coquimbo'
$ws.Cells.Item(31,4).Value = 44242
$ws.Cells.Item(31,5).Value = 5
$ws.Cells.Item(31,6).Value = 'Fruta'
$ws.Cells.Item(31,7).Value = 100104
$ws.Cells.Item(31,8).Value = 'Frutos de pepita'
$ws.Cells.Item(31,9).Value = 100104001
$ws.Cells.Item(31,10).Value = 'Granada'
$ws.Cells.Item(31,11).Value = 'Wonderfull'
$ws.Cells.Item(31,12).Value = 'Primera'
$ws.Cells.Item(31,13).Value = 45
$ws.Cells.Item(31,14).Value = 12000
$ws.Cells.Item(31,15).Value = 12000
$ws.Cells.Item(31,16).Value = 12000
$ws.Cells.Item(31,17).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(31,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(31,19).Value = 800
$ws.Cells.Item(31,20).Value = 15
$ws.Cells.Item(32,1).Value = 3
$ws.Cells.Item(32,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(32,3).Value = 'Coquimbo'
$ws.Cells.Item(32,4).Value = 45040
$ws.Cells.Item(32,5).Value = 5
$ws.Cells.Item(32,6).Value = 'Fruta'
$ws.Cells.Item(32,7).Value = 100104
$ws.Cells.Item(32,8).Value = 'Frutos de pepita'
$ws.Cells.Item(32,9).Value = 100104001
$ws.Cells.Item(32,10).Value = 'Granada'
$ws.Cells.Item(32,11).Value = 'Wonderfull'
$ws.Cells.Item(32,12).Value = 'Especial'
$ws.Cells.Item(32,13).Value = 65
$ws.Cells.Item(32,14).Value = 17000
$ws.Cells.Item(32,15).Value = 17000
$ws.Cells.Item(32,16).Value = 17000
$ws.Cells.Item(32,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(32,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(32,19).Value = 1214
$ws.Cells.Item(32,20).Value = 14
$ws.Cells.Item(33,1).Value = 3
$ws.Cells.Item(33,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(33,3).Value = 'Coquimbo'
$ws.Cells.Item(33,4).Value = 45040
$ws.Cells.Item(33,5).Value = 5
$ws.Cells.Item(33,6).Value = 'Fruta'
$ws.Cells.Item(33,7).Value = 100104
$ws.Cells.Item(33,8).Value = 'Frutos de pepita'
$ws.Cells.Item(33,9).Value = 100104001
$ws.Cells.Item(33,10).Value = 'Granada'
$ws.Cells.Item(33,11).Value = 'Wonderfull'
$ws.Cells.Item(33,12).Value = 'Primera'
$ws.Cells.Item(33,13).Value = 60
$ws.Cells.Item(33,14).Value = 14000
$ws.Cells.Item(33,15).Value = 14000
$ws.Cells.Item(33,16).Value = 14000
$ws.Cells.Item(33,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(33,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(33,19).Value = 1000
$ws.Cells.Item(33,20).Value = 14
$ws.Cells.Item(34,1).Value = 3
$ws.Cells.Item(34,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(34,3).Value = 'Coquimbo'
$ws.Cells.Item(34,4).Value = 44616
$ws.Cells.Item(34,5).Value = 5
$ws.Cells.Item(34,6).Value = 'Fruta'
$ws.Cells.Item(34,7).Value = 100104
$ws.Cells.Item(34,8).Value = 'Frutos de pepita'
$ws.Cells.Item(34,9).Value = 100104001
$ws.Cells.Item(34,10).Value = 'Granada'
$ws.Cells.Item(34,11).Value = 'Wonderfull'
$ws.Cells.Item(34,12).Value = 'Primera'
$ws.Cells.Item(34,13).Value = 70
$ws.Cells.Item(34,14).Value = 14000
$ws.Cells.Item(34,15).Value = 14000
$ws.Cells.Item(34,16).Value = 14000
$ws.Cells.Item(34,17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(34,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(34,19).Value = 1000
$ws.Cells.Item(34,20).Value = 14
$ws.Cells.Item(35,1).Value = 3
$ws.Cells.Item(35,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(35,3).Value = 'Coquimbo'
$ws.Cells.Item(35,4).Value = 44323
$ws.Cells.Item(35,5).Value = 5
$ws.Cells.Item(35,6).Value = 'Fruta'
$ws.Cells.Item(35,7).Value = 100104
$ws.Cells.Item(35,8).Value = 'Frutos de pepita'
$ws.Cells.Item(35,9).Value = 100104001
$ws.Cells.Item(35,10).Value = 'Granada'
$ws.Cells.Item(35,11).Value = 'Wonderfull'
$ws.Cells.Item(35,12).Value = 'Primera'
$ws.Cells.Item(35,13).Value = 60
$ws.Cells.Item(35,14).Value = 14000
$ws.Cells.Item(35,15).Value = 14000
$ws.Cells.Item(35,16).Value = 14000
$ws.Cells.Item(35,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(35,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(35,19).Value = 1000
$ws.Cells.Item(35,20).Value = 14
$ws.Cells.Item(36,1).Value = 3
$ws.Cells.Item(36,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(36,3).Value = 'Coquimbo'
$ws.Cells.Item(36,4).Value = 44315
$ws.Cells.Item(36,5).Value = 5
$ws.Cells.Item(36,6).Value = 'Fruta'
$ws.Cells.Item(36,7).Value = 100104
$ws.Cells.Item(36,8).Value = 'Frutos de pepita'
$ws.Cells.Item(36,9).Value = 100104001
$ws.Cells.Item(36,10).Value = 'Granada'
$ws.Cells.Item(36,11).Value = 'Wonderfull'
$ws.Cells.Item(36,12).Value = 'Primera'
$ws.Cells.Item(36,13).Value = 65
$ws.Cells.Item(36,14).Value = 14000
$ws.Cells.Item(36,15).Value = 14000
$ws.Cells.Item(36,16).Value = 14000
$ws.Cells.Item(36,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(36,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(36,19).Value = 1000
$ws.Cells.Item(36,20).Value = 14
$ws.Cells.Item(37,1).Value = 3
$ws.Cells.Item(37,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(37,3).Value = 'Coquimbo'
$ws.Cells.Item(37,4).Value = 44592
$ws.Cells.Item(37,5).Value = 5
$ws.Cells.Item(37,6).Value = 'Fruta'
$ws.Cells.Item(37,7).Value = 100104
$ws.Cells.Item(37,8).Value = 'Frutos de pepita'
$ws.Cells.Item(37,9).Value = 100104001
$ws.Cells.Item(37,10).Value = 'Granada'
$ws.Cells.Item(37,11).Value = 'Wonderfull'
$ws.Cells.Item(37,12).Value = 'Primera'
$ws.Cells.Item(37,13).Value = 54
$ws.Cells.Item(37,14).Value = 20000
$ws.Cells.Item(37,15).Value = 20000
$ws.Cells.Item(37,16).Value = 20000
$ws.Cells.Item(37,17).Value = '$/caja 15 kilos empedrada'
$ws.Cells.Item(37,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(37,19).Value = 1333
$ws.Cells.Item(37,20).Value = 15
$ws.Cells.Item(38,1).Value = 3
$ws.Cells.Item(38,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(38,3).Value = 'Coquimbo'
$ws.Cells.Item(38,4).Value = 44271
$ws.Cells.Item(38,5).Value = 5
$ws.Cells.Item(38,6).Value = 'Fruta'
$ws.Cells.Item(38,7).Value = 100104
$ws.Cells.Item(38,8).Value = 'Frutos de pepita'
$ws.Cells.Item(38,9).Value = 100104001
$ws.Cells.Item(38,10).Value = 'Granada'
$ws.Cells.Item(38,11).Value = 'Wonderfull'
$ws.Cells.Item(38,12).Value = 'Primera'
$ws.Cells.Item(38,13).Value = 50
$ws.Cells.Item(38,14).Value = 12000
$ws.Cells.Item(38,15).Value = 12000
$ws.Cells.Item(38,16).Value = 12000
$ws.Cells.Item(38,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(38,18).Value = 'Provincia del Elquí'
$ws.Cells.Item(38,19).Value = 857
$ws.Cells.Item(38,20).Value = 14
$ws.Cells.Item(39,1).Value = 3
$ws.Cells.Item(39,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(39,3).Value = 'Coquimbo'
$ws.Cells.Item(39,4).Value = 44322
$ws.Cells.Item(39,5).Value = 5
$ws.Cells.Item(39,6).Value = 'Fruta'
$ws.Cells.Item(39,7).Value = 100104
$ws.Cells.Item(39,8).Value = 'Frutos de pepita'
$ws.Cells.Item(39,9).Value = 100104001
$ws.Cells.Item(39,10).Value = 'Granada'
$ws.Cells.Item(39,11).Value = 'Wonderfull'
$ws.Cells.Item(39,12).Value = 'Primera'
$ws.Cells.Item(39,13).Value = 50
$ws.Cells.Item(39,14).Value = 14000
$ws.Cells.Item(39,15).Value = 14000
$ws.Cells.Item(39,16).Value = 14000
$ws.Cells.Item(39,17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item(39,18).Value = 'Provincia de Limarí'
$ws.Cells.Item(39,19).Value = 1000
$ws.Cells.Item(39,20).Value = 14
